# "Fixed up sankey and tableau for time"
# The existing "2015" sheet (Uber day-of-week totals) becomes "2014" with new
# 2014 data, and a fresh "2015" sheet is added (after it) holding the
# original data that used to live on the single "2015" sheet.

$wb = $excel.ActiveWorkbook

# --- Step 1: repurpose the existing sheet as "2014" --------------------
$sheet2014 = $wb.Worksheets.Item(1)
$sheet2014.Name = "2014"

# --- Step 2: add a new sheet right after it, named "2015" --------------
$sheet2015 = $wb.Worksheets.Add($null, $sheet2014)
$sheet2015.Name = "2015"

# --- Step 3: move the original ("2015") day_of_week data onto the new --
#             "2015" sheet, unchanged.
$orig = @(
    @("Sat", 2414563),
    @("Fri", 2282571),
    @("Thu", 2159598),
    @("Sun", 1952782),
    @("Wed", 1893811),
    @("Tue", 1872902),
    @("Mon", 1694252)
)

$sheet2015.Range("B1").Value = "day_of_week"
for ($i = 0; $i -lt $orig.Length; $i++) {
    $row = $i + 2
    $sheet2015.Cells.Item($row, 1).Value = $orig[$i][0]
    $sheet2015.Cells.Item($row, 2).Value = $orig[$i][1]
}

# --- Step 4: overwrite the "2014" sheet with the new 2014 data ---------
$data2014 = @(
    @("Thu", 755145),
    @("Fri", 741139),
    @("Wed", 696488),
    @("Tue", 663789),
    @("Sat", 646114),
    @("Mon", 541472),
    @("Sun", 490180)
)

$sheet2014.Range("B1").Value = "day_of_week"
for ($i = 0; $i -lt $data2014.Length; $i++) {
    $row = $i + 2
    $sheet2014.Cells.Item($row, 1).Value = $data2014[$i][0]
    $sheet2014.Cells.Item($row, 2).Value = $data2014[$i][1]
}

# --- Step 5: apply the header/label style (bold, thin box border, ------
#             centered horizontally, top-aligned vertically) to the
#             label cells on both sheets, matching the original format.
#             (Union ranges like "B1,A2:A8" only format their first area in
#             this host, so the two blocks are styled individually.)
foreach ($ws in @($sheet2014, $sheet2015)) {
    foreach ($rng in @($ws.Range("B1"), $ws.Range("A2:A8"))) {
        $rng.Font.Bold = $true
        $rng.Borders.LineStyle = 1
        $rng.HorizontalAlignment = -4108
        $rng.VerticalAlignment = -4160
    }
}

# --- Step 6: make "2014" the active sheet/tab (activeTab=0) ------------
$sheet2014.Activate()
